$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new CKJ51 row (row 27), inserted logically between CKJ46 and CKJ52
$ws.Range("A27").Value = "CKJ51"

# Refresh call-volume figures for every existing row plus the new CKJ51 row
$data = @{
    2 = @(0, 2, 3, 1, 6, 1, 6, 0, 3)
    3 = @(6, 2, 5, 5, 2, 5, 5, 5, 2)
    4 = @(7, 1, 3, 7, 5, 7, 7, 9, 9)
    5 = @(15, 6, 13, 5, 7, 13, 12, 5, 9)
    6 = @(8, 3, 6, 1, 1, 0, 10, 1, 5)
    7 = @(9, 12, 6, 13, 15, 7, 4, 0, 0)
    8 = @(9, 0, 3, 4, 3, 4, 8, 0, 7)
    9 = @(9, 4, 4, 10, 13, 10, 9, 8, 1)
    10 = @(2, 1, 2, 5, 8, 8, 8, 4, 3)
    11 = @(1, 1, 3, 4, 5, 5, 6, 0, 0)
    12 = @(1, 3, 4, 0, 2, 9, 7, 1, 4)
    13 = @(5, 2, 4, 4, 1, 6, 3, 1, 1)
    14 = @(12, 5, 0, 10, 12, 10, 10, 0, 0)
    15 = @(10, 5, 5, 9, 6, 12, 12, 5, 5)
    16 = @(5, 5, 3, 3, 4, 6, 11, 0, 0)
    17 = @(9, 2, 10, 0, 7, 3, 3, 1, 5)
    18 = @(1, 1, 2, 1, 5, 4, 15, 6, 15)
    19 = @(9, 11, 3, 15, 0, 8, 4, 0, 5)
    20 = @(9, 3, 3, 3, 4, 8, 12, 1, 5)
    21 = @(10, 5, 8, 2, 3, 19, 9, 1, 13)
    22 = @(12, 7, 4, 9, 6, 11, 9, 9, 4)
    23 = @(13, 14, 9, 11, 16, 16, 13, 4, 5)
    24 = @(0, 0, 0, 5, 2, 3, 12, 0, 0)
    25 = @(13, 7, 11, 8, 2, 14, 1, 9, 1)
    26 = @(6, 5, 2, 1, 3, 7, 6, 1, 5)
    27 = @(8, 4, 5, 5, 19, 12, 7, 11, 8)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 2).Value = $values[$i]
    }
}